$d = $word.ActiveDocument

# The document mentions "tretåig hackspett" (three-toed woodpecker) as the
# single prioritized species; the edit adds "talltita" (willow tit) as a
# second species throughout, and updates surrounding singular->plural
# Swedish grammar (art -> arter, denna -> dessa, arten -> arterna), plus
# bumps the trailing date by one day.

# 1. "tretåig hackspett" -> "talltita och tretåig hackspett" (both occurrences
#    in the document: the main body paragraph and the bulleted list item)
$d.Content.Find.Execute("tretåig hackspett", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "talltita och tretåig hackspett", 2)

# 2. "Detta är en prioriterad art" -> "Dessa är prioriterade arter"
$d.Content.Find.Execute("Detta är en prioriterad art", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Dessa är prioriterade arter", 2)

# 3. " denna art" -> " dessa arter"
$d.Content.Find.Execute(" denna art", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " dessa arter", 2)

# 4. " arten" -> " arterna"
$d.Content.Find.Execute(" arten", $false, $false, $false, $false, $false, `
                         $true, 1, $false, " arterna", 2)

# 5. Update the signature date "2026-02-20" -> "2026-02-21"
$d.Content.Find.Execute("2026-02-20", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-02-21", 2)
